$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: now holds the record previously in row 8 (Stiftgelélav / Collema furfuraceum),
# but with its own new Id/Taxonsorteringsordning/coordinates/time/comment/observer values.
$ws.Range("A7").Value = 131108352
$ws.Range("B7").Value = 80216
$ws.Range("E7").Value = 388
$ws.Range("F7").Value = 'Stiftgelélav'
$ws.Range("G7").Value = 'Collema furfuraceum'
$ws.Range("H7").Value = '(Arnold) Du Rietz'
$ws.Range("I7").NumberFormat = "@"
$ws.Range("I7").Value = '1'
$ws.Range("J7").Value = 'bålar'
$ws.Range("P7").Value = 'S Svartmyran, Mpd'
$ws.Range("Q7").Value = 616863
$ws.Range("R7").Value = 6934788
$ws.Range("X7").Value = '2025_0758'
$ws.Range("Z7").Value = '14:47'
$ws.Range("AB7").Value = '14:47'
$ws.Range("AC7").Value = 'Asp'
$ws.Range("AX7").Value = 'Måns Svensson'

# Row 8: now holds the record previously in row 7 (Reliktbock / Nothorhina muricata).
$ws.Range("A8").Value = 131106436
$ws.Range("B8").Value = 5493
$ws.Range("E8").Value = 101410
$ws.Range("F8").Value = 'Reliktbock'
$ws.Range("G8").Value = 'Nothorhina muricata'
$ws.Range("H8").Value = '(Dalman, 1817)'
$ws.Range("I8").NumberFormat = "@"
$ws.Range("I8").Value = '2'
$ws.Range("J8").Value = 'ex.'
$ws.Range("P8").Value = 'Svartmyran, Mpd'
$ws.Range("Q8").Value = 616762
$ws.Range("R8").Value = 6934714
$ws.Range("X8").Value = '2025_0743'
$ws.Range("Z8").Value = '11:39'
$ws.Range("AB8").Value = '11:39'
$ws.Range("AC8").Value = 'Två kläckhål'
$ws.Range("AX8").Value = 'David Isaksson'

# Rows 9-11: Taxonsorteringsordning bumped by 1.
$ws.Range("B9").Value = 80254
$ws.Range("B10").Value = 80379
$ws.Range("B11").Value = 80350
